$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.385.26'
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.485.51'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.95'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.68'
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.509'
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.486.60'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("E10").Value = '  +2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.336'
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.89'
$ws.Range("E13").Value = '  +3.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '69.249.99'
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.88'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.481.64'
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.88'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.58'
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.11'
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.83'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("E23").Value = '  +6.40%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.65'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.75'
$ws.Range("E26").Value = '  +1.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.611.20'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.34'
$ws.Range("E28").Value = '  +2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0834'
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.27'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '438.85'
$ws.Range("E32").Value = '  +3.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.16'
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  -1.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.65'
$ws.Range("E36").Value = '  +0.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.07'
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.110'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.97'
$ws.Range("E40").Value = '  +1.28%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.49'
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.304'
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '37.65'
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.50'
$ws.Range("E44").Value = '  -0.53%  '
$ws.Range("E45").Value = '  +4.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.11'
$ws.Range("E46").Value = '  +3.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.44'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.38'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0719'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.490'
$ws.Range("E50").Value = '  +2.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.567'
$ws.Range("E51").Value = '  +2.47%  '
